# Demo Verification Script / Demo-QA CF fixes:
# Update the "Date" (B2) timestamp on the relevant test-result sheets to the
# latest verification run, and flip CMCAutopayPS's Result (A2) from Fail to
# Pass now that the underlying issue is fixed.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "PayNowNoCFPC";              Date = "Wed Jun 25 00:51:14 IST 2025" },
    @{ Sheet = "PayNowNoCFPS";              Date = "Wed Jun 25 00:52:03 IST 2025" },
    @{ Sheet = "PayNowSCFPC";               Date = "Wed Jun 25 00:54:16 IST 2025" },
    @{ Sheet = "PayNowSCFPS";               Date = "Wed Jun 25 00:55:14 IST 2025" },
    @{ Sheet = "PayNowDCFPC";               Date = "Wed Jun 25 00:44:02 IST 2025" },
    @{ Sheet = "PayNowDCFPS";               Date = "Wed Jun 25 00:44:52 IST 2025" },
    @{ Sheet = "SCFPSVerbiage";             Date = "Wed Jun 25 01:00:05 IST 2025" },
    @{ Sheet = "SCFPCVerbiage";             Date = "Wed Jun 25 00:59:27 IST 2025" },
    @{ Sheet = "DCFPSVerbiage";             Date = "Wed Jun 25 00:57:44 IST 2025" },
    @{ Sheet = "DCFPCVerbiage";             Date = "Wed Jun 25 00:57:24 IST 2025" },
    @{ Sheet = "CMCAutopayPS";              Date = "Wed Jun 25 00:33:53 IST 2025"; Result = "Pass" },
    @{ Sheet = "CMCAutoPayPC";              Date = "Wed Jun 25 00:32:18 IST 2025" },
    @{ Sheet = "NoModifyAmountPC";          Date = "Wed Jun 25 00:37:04 IST 2025" },
    @{ Sheet = "NoModifyAmountPS";          Date = "Wed Jun 25 00:38:05 IST 2025" },
    @{ Sheet = "NoModifyBillingAddressPC";  Date = "Wed Jun 25 00:40:37 IST 2025" },
    @{ Sheet = "NoModifyBillingAddressPS";  Date = "Wed Jun 25 00:41:33 IST 2025" }
)

foreach ($update in $updates) {
    $ws = $wb.Worksheets.Item($update.Sheet)
    if ($update.ContainsKey("Result")) {
        $ws.Range("A2").Value = $update.Result
    }
    $ws.Range("B2").Value = $update.Date
}
